# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these Price cells to remain plain text so numeric-looking values
# ("1.004", "0.9650", "0.000009882", ...) keep their exact original formatting
# instead of being auto-converted to numbers by Excel.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated values
$ws.Range("D2").Value = "19.757.61"
$ws.Range("D3").Value = "1.390.64"
$ws.Range("E3").Value = "  -9.34%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.37%  "
$ws.Range("D5").Value = "1.004"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D6").Value = "268.62"
$ws.Range("E6").Value = "  -6.83%  "
$ws.Range("D7").Value = "0.3643"
$ws.Range("E7").Value = "  -7.42%  "
$ws.Range("D8").Value = "0.3035"
$ws.Range("E8").Value = "  -3.90%  "
$ws.Range("D9").Value = "38.43"
$ws.Range("E9").Value = "  -9.43%  "
$ws.Range("D10").Value = "0.06394"
$ws.Range("E10").Value = "  -10.71%  "
$ws.Range("D11").Value = "0.9650"
$ws.Range("E11").Value = "  -7.67%  "
$ws.Range("D12").Value = "1.005"
$ws.Range("E12").Value = "  +0.44%  "
$ws.Range("D13").Value = "5.278"
$ws.Range("E13").Value = "  -6.62%  "
$ws.Range("D14").Value = "6.036"
$ws.Range("E14").Value = "  -8.31%  "
$ws.Range("D15").Value = "1.393.31"
$ws.Range("D16").Value = "16.36"
$ws.Range("E16").Value = "  -11.78%  "
$ws.Range("D17").Value = "0.000009882"
$ws.Range("E17").Value = "  -9.22%  "
$ws.Range("D18").Value = "0.05626"
$ws.Range("E18").Value = "  -14.82%  "
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("D20").Value = "70.09"
$ws.Range("E20").Value = "  -16.07%  "
$ws.Range("D21").Value = "5.497"
$ws.Range("E21").Value = "  -10.03%  "
$ws.Range("D22").Value = "14.18"
$ws.Range("E22").Value = "  -8.01%  "
$ws.Range("E23").Value = "  -2.79%  "
$ws.Range("E24").Value = "  -4.70%  "
$ws.Range("D25").Value = "19.750.53"
$ws.Range("E25").Value = "  -8.76%  "
$ws.Range("D26").Value = "2.142"
$ws.Range("E26").Value = "  -8.86%  "
$ws.Range("D27").Value = "135.65"
$ws.Range("E27").Value = "  -8.46%  "
$ws.Range("D28").Value = "16.51"
$ws.Range("E28").Value = "  -9.83%  "
$ws.Range("D29").Value = "1.546.40"
$ws.Range("E29").Value = "  -9.59%  "
$ws.Range("D30").Value = "107.42"
$ws.Range("E30").Value = "  -8.24%  "
$ws.Range("D31").Value = "3.841"
$ws.Range("E31").Value = "  -20.65%  "
$ws.Range("D32").Value = "5.230"
$ws.Range("E32").Value = "  -11.90%  "
$ws.Range("D33").Value = "0.7851"
$ws.Range("E33").Value = "  -16.96%  "
$ws.Range("D34").Value = "0.07581"
$ws.Range("E34").Value = "  -6.92%  "
$ws.Range("D35").Value = "8.187"
$ws.Range("E35").Value = "  -3.47%  "
$ws.Range("E36").Value = "  +0.48%  "
$ws.Range("D37").Value = "4.701"
$ws.Range("E37").Value = "  -8.33%  "
$ws.Range("D38").Value = "0.05537"
$ws.Range("E38").Value = "  -7.71%  "
$ws.Range("D39").Value = "0.02006"
$ws.Range("E39").Value = "  -8.98%  "
$ws.Range("D40").Value = "0.1861"
$ws.Range("E40").Value = "  -7.55%  "
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "9.993"
$ws.Range("E41").Value = "  -9.09%  "
$ws.Range("B42").Value = "WEMIXTOKEN"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "1.282"
$ws.Range("E42").Value = "  -11.20%  "
$ws.Range("D43").Value = "1.028"
$ws.Range("E43").Value = "  -12.31%  "
$ws.Range("D44").Value = "3.465"
$ws.Range("E44").Value = "  -6.46%  "
$ws.Range("E45").Value = "  -10.43%  "
$ws.Range("D46").Value = "11.75"
$ws.Range("E46").Value = "  -9.98%  "
$ws.Range("D47").Value = "0.4941"
$ws.Range("E47").Value = "  -9.79%  "
$ws.Range("D48").Value = "108.09"
$ws.Range("E48").Value = "  -6.92%  "
$ws.Range("D49").Value = "1.708"
$ws.Range("E49").Value = "  -8.74%  "
$ws.Range("D50").Value = "1.005"
$ws.Range("E50").Value = "  +0.47%  "
$ws.Range("D51").Value = "1.029"
$ws.Range("E51").Value = "  -11.57%  "
